$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibition) - update "想去人数" (interested-count) values
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 2728
$wsExpo.Range("F5").Value  = 938
$wsExpo.Range("F7").Value  = 2343
$wsExpo.Range("F8").Value  = 1845
$wsExpo.Range("F10").Value = 62
$wsExpo.Range("F11").Value = 2497
$wsExpo.Range("F17").Value = 119
$wsExpo.Range("F18").Value = 9293
$wsExpo.Range("F20").Value = 7216
$wsExpo.Range("F21").Value = 11776
$wsExpo.Range("F24").Value = 236
$wsExpo.Range("F25").Value = 363
$wsExpo.Range("F26").Value = 563
$wsExpo.Range("F27").Value = 2631
$wsExpo.Range("F29").Value = 200
$wsExpo.Range("F30").Value = 2572
$wsExpo.Range("F31").Value = 740
$wsExpo.Range("F33").Value = 4523
$wsExpo.Range("F34").Value = 946
$wsExpo.Range("F35").Value = 356
$wsExpo.Range("F36").Value = 44
$wsExpo.Range("F37").Value = 538

# ---------------------------------------------------------------------------
# Sheet "演出" (Performance) - value update + a brand-new event inserted
# as row 18 (everything below shifts down by one row)
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 70

# Insert a blank row at position 18; rows 18-25 move down to 19-26.
$wsShow.Rows(18).Insert()

# The freshly-inserted row's index cell (A18) comes back blank/unstyled.
# Restore its look-and-feel by copying a neighbouring, correctly-styled
# index cell onto it (value gets fixed up below).
$wsShow.Cells.Item(17, 1).Copy($wsShow.Cells.Item(18, 1))

# Populate the new row with the "蕨野友也粉丝见面会" event.
# The leading "'" keeps the date-shaped text as a literal string instead of
# letting it auto-convert to a date serial.
$wsShow.Cells.Item(18, 2).Value = "'2024-12-21"
$wsShow.Cells.Item(18, 3).Value = "杭州·蕨野友也粉丝见面会"
$wsShow.Cells.Item(18, 4).Value = "教工路198号浙商大创业园A幢3楼 杭州子墨汇演中心"
$wsShow.Cells.Item(18, 5).Value = "2024.12.21 12:30-12.21 17:30"
$wsShow.Cells.Item(18, 6).Value = 0
$wsShow.Cells.Item(18, 7).Value = 380
$wsShow.Cells.Item(18, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93455"
$wsShow.Cells.Item(18, 9).Value = "//i2.hdslb.com/bfs/openplatform/202410/2t3vGbf21728964006170.jpeg"

# Column A is a plain row-index ("row number - 1"); it is NOT one of the
# fields that travels with the rest of a row's data, so after the insert
# every index cell from row 18 down must be reset to (row - 1), including
# the brand-new row 26 created by the shift.
for ($r = 18; $r -le 26; $r++) {
    $wsShow.Cells.Item($r, 1).Value = $r - 1
}

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local Life)
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F4").Value = 164

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All Types) - update "想去人数" values
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 2728
$wsAll.Range("F6").Value  = 70
$wsAll.Range("F7").Value  = 938
$wsAll.Range("F10").Value = 2343
$wsAll.Range("F12").Value = 1845
$wsAll.Range("F15").Value = 2497
$wsAll.Range("F22").Value = 119
$wsAll.Range("F23").Value = 9293
$wsAll.Range("F25").Value = 7216
$wsAll.Range("F26").Value = 11777
$wsAll.Range("F29").Value = 236
$wsAll.Range("F30").Value = 363
$wsAll.Range("F32").Value = 563
$wsAll.Range("F34").Value = 2631
$wsAll.Range("F38").Value = 200
$wsAll.Range("F40").Value = 4523
$wsAll.Range("F45").Value = 538
